$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update credentials shown in the sheet
$ws.Range("A2").Value = "anisa@gmail.com"
$ws.Range("B2").Value = "xyz123"

# Update the selected cell to match the new active selection
$ws.Range("C3").Select()
